$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend column widths for new columns J:L (width 35) ---
$ws.Columns.Item(10).ColumnWidth = 34.16
$ws.Columns.Item(11).ColumnWidth = 34.16
$ws.Columns.Item(12).ColumnWidth = 34.16

# --- 2. Copy formatting (style) for new header cells J1:L1 from existing header cell I1 ---
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:L1").PasteSpecial(-4122) | Out-Null

# --- 3. Copy formatting for new row headers A7:A12 from existing row header A6 ---
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7:A12").PasteSpecial(-4122) | Out-Null

# --- 4. Copy formatting for new data cells (rows 7-12 cols B-I, and all rows cols J-L) from existing data cell B2 ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B7:L12").PasteSpecial(-4122) | Out-Null
$ws.Range("J2:L6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 5. Set header row values ---
$ws.Range("B1").Value = "css3"
$ws.Range("C1").Value = "external"
$ws.Range("D1").Value = "google-workspace"
$ws.Range("E1").Value = "html5"
$ws.Range("F1").Value = "internet"
$ws.Range("G1").Value = "javascript"
$ws.Range("H1").Value = "jboss"
$ws.Range("I1").Value = "keycloak"
$ws.Range("J1").Value = "linux"
$ws.Range("K1").Value = "oauth"
$ws.Range("L1").Value = "webrtc"

# --- 6. Set row headers and data grid values ---
$ws.Range("A2").Value = "Auth Server"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "X"
$ws.Range("I2").Value = "X"
$ws.Range("J2").Value = "X"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

$ws.Range("A3").Value = "Google Meet Client"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "X"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "X"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""

$ws.Range("A4").Value = "Google Meet Server"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""

$ws.Range("A5").Value = "Web Interface"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""

$ws.Range("A6").Value = "User Authentication Access"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""

$ws.Range("A7").Value = "Load Balancer"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""

$ws.Range("A8").Value = "Google Meet Traffic"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""

$ws.Range("A9").Value = "User Account Information"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "X"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = "X"
$ws.Range("L9").Value = ""

$ws.Range("A10").Value = "Video Meeting Data"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = "X"

$ws.Range("A11").Value = "Public Internet"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "X"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "X"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""

$ws.Range("A12").Value = "Google Cloud Platform"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""

